# contratos-7-2017.xlsx -- fix formatting glitches introduced when scraping data:
#  1) four provider names used a comma where a period was intended
#  2) the "Importe" column (H) held amounts as text using the Spanish
#     "1.234,56" convention; re-save them as plain "1234.56" text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Provider-name punctuation fixes (comma -> period) ---
$ws.Range("E80").Value = 'FERNANDEZ. MARIO HUGO'
$ws.Range("E148").Value = 'DODERA. JORGE ABELARDO'
$ws.Range("E166").Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range("E225").Value = 'LUGRIN. RICARDO JOSE'
$ws.Range("F225").Value = 'LUGRIN. RICARDO JOSE'

# --- 2) Re-format the Importe (H) text amounts ---
# These cells hold text, not real numbers. A plain Value= assignment of a
# digit string would make Excel re-interpret it as a number and silently
# drop the trailing ".00", so prefix with an apostrophe to keep it text.
$ws.Range("H2").Value = "'2520.00"
$ws.Range("H3").Value = "'930.00"
$ws.Range("H4").Value = "'9445.00"
$ws.Range("H5").Value = "'35900.00"
$ws.Range("H6").Value = "'10500.00"
$ws.Range("H7").Value = "'1385000.00"
$ws.Range("H8").Value = "'429500.00"
$ws.Range("H9").Value = "'950000.00"
$ws.Range("H10").Value = "'128000.00"
$ws.Range("H11").Value = "'469000.00"
$ws.Range("H12").Value = "'388727.60"
$ws.Range("H13").Value = "'27225.00"
$ws.Range("H14").Value = "'698000.00"
$ws.Range("H15").Value = "'5750.02"
$ws.Range("H16").Value = "'240.00"
$ws.Range("H17").Value = "'859678.20"
$ws.Range("H18").Value = "'516197.88"
$ws.Range("H19").Value = "'27956.00"
$ws.Range("H20").Value = "'87436.40"
$ws.Range("H21").Value = "'73757.50"
$ws.Range("H22").Value = "'59639.84"
$ws.Range("H23").Value = "'399.00"
$ws.Range("H24").Value = "'33355.12"
$ws.Range("H25").Value = "'19050.00"
$ws.Range("H26").Value = "'4510.00"
$ws.Range("H27").Value = "'7500.00"
$ws.Range("H28").Value = "'7500.00"
$ws.Range("H29").Value = "'505.31"
$ws.Range("H30").Value = "'43400.00"
$ws.Range("H31").Value = "'9680.00"
$ws.Range("H32").Value = "'22442.10"
$ws.Range("H33").Value = "'1840.00"
$ws.Range("H34").Value = "'24610.00"
$ws.Range("H35").Value = "'355810.00"
$ws.Range("H36").Value = "'1077.00"
$ws.Range("H37").Value = "'1760.00"
$ws.Range("H38").Value = "'6458.89"
$ws.Range("H39").Value = "'800.00"
$ws.Range("H40").Value = "'388.56"
$ws.Range("H41").Value = "'4738.00"
$ws.Range("H42").Value = "'10197.75"
$ws.Range("H43").Value = "'1299.00"
$ws.Range("H44").Value = "'153915.54"
$ws.Range("H45").Value = "'21888.00"
$ws.Range("H46").Value = "'2200.00"
$ws.Range("H47").Value = "'2683.00"
$ws.Range("H48").Value = "'8734.28"
$ws.Range("H49").Value = "'11230.77"
$ws.Range("H50").Value = "'1065.00"
$ws.Range("H51").Value = "'16203.43"
$ws.Range("H52").Value = "'6101.80"
$ws.Range("H53").Value = "'7631.10"
$ws.Range("H54").Value = "'38102.55"
$ws.Range("H55").Value = "'3060.00"
$ws.Range("H56").Value = "'20052.18"
$ws.Range("H57").Value = "'330.00"
$ws.Range("H58").Value = "'1321.84"
$ws.Range("H59").Value = "'471.88"
$ws.Range("H60").Value = "'4810.00"
$ws.Range("H61").Value = "'42000.00"
$ws.Range("H62").Value = "'3270.00"
$ws.Range("H63").Value = "'90256.00"
$ws.Range("H64").Value = "'11500.00"
$ws.Range("H65").Value = "'4260.00"
$ws.Range("H66").Value = "'5620.00"
$ws.Range("H67").Value = "'596.16"
$ws.Range("H68").Value = "'433.16"
$ws.Range("H69").Value = "'35500.00"
$ws.Range("H70").Value = "'2350.00"
$ws.Range("H71").Value = "'980.00"
$ws.Range("H72").Value = "'3036.15"
$ws.Range("H73").Value = "'1500.00"
$ws.Range("H74").Value = "'169576.00"
$ws.Range("H75").Value = "'656.00"
$ws.Range("H76").Value = "'117227.00"
$ws.Range("H77").Value = "'24015.00"
$ws.Range("H78").Value = "'66132.00"
$ws.Range("H79").Value = "'800.00"
$ws.Range("H80").Value = "'8100.00"
$ws.Range("H81").Value = "'62190.00"
$ws.Range("H82").Value = "'482.00"
$ws.Range("H83").Value = "'12.00"
$ws.Range("H84").Value = "'10065.00"
$ws.Range("H85").Value = "'558358.95"
$ws.Range("H86").Value = "'16851.11"
$ws.Range("H87").Value = "'8.00"
$ws.Range("H88").Value = "'220.70"
$ws.Range("H89").Value = "'78.18"
$ws.Range("H90").Value = "'31488.69"
$ws.Range("H91").Value = "'203.00"
$ws.Range("H92").Value = "'280.00"
$ws.Range("H93").Value = "'60.00"
$ws.Range("H94").Value = "'5213.00"
$ws.Range("H95").Value = "'30742.12"
$ws.Range("H96").Value = "'5240.00"
$ws.Range("H97").Value = "'919.30"
$ws.Range("H98").Value = "'53.80"
$ws.Range("H99").Value = "'3013.80"
$ws.Range("H100").Value = "'794.00"
$ws.Range("H101").Value = "'24583.00"
$ws.Range("H102").Value = "'12265.00"
$ws.Range("H103").Value = "'8910.00"
$ws.Range("H104").Value = "'94350.00"
$ws.Range("H105").Value = "'3500.00"
$ws.Range("H106").Value = "'25000.00"
$ws.Range("H107").Value = "'2500.00"
$ws.Range("H108").Value = "'7200.00"
$ws.Range("H109").Value = "'162500.00"
$ws.Range("H110").Value = "'600.00"
$ws.Range("H111").Value = "'6134.33"
$ws.Range("H112").Value = "'1761.10"
$ws.Range("H113").Value = "'8500.00"
$ws.Range("H114").Value = "'1800.00"
$ws.Range("H115").Value = "'435.00"
$ws.Range("H116").Value = "'15205.00"
$ws.Range("H117").Value = "'20560.00"
$ws.Range("H118").Value = "'10000.00"
$ws.Range("H119").Value = "'6000.00"
$ws.Range("H120").Value = "'14000.00"
$ws.Range("H121").Value = "'6000.00"
$ws.Range("H122").Value = "'25116.58"
$ws.Range("H123").Value = "'4500.00"
$ws.Range("H124").Value = "'3000.00"
$ws.Range("H125").Value = "'4999.12"
$ws.Range("H126").Value = "'4000.00"
$ws.Range("H127").Value = "'3500.00"
$ws.Range("H128").Value = "'3000.00"
$ws.Range("H129").Value = "'3000.00"
$ws.Range("H130").Value = "'4500.00"
$ws.Range("H131").Value = "'20000.00"
$ws.Range("H132").Value = "'3500.00"
$ws.Range("H133").Value = "'3000.00"
$ws.Range("H134").Value = "'4500.00"
$ws.Range("H135").Value = "'46907.25"
$ws.Range("H136").Value = "'10500.00"
$ws.Range("H137").Value = "'6000.00"
$ws.Range("H138").Value = "'4500.00"
$ws.Range("H139").Value = "'6000.00"
$ws.Range("H140").Value = "'4000.00"
$ws.Range("H141").Value = "'15000.00"
$ws.Range("H142").Value = "'15000.00"
$ws.Range("H143").Value = "'24300.00"
$ws.Range("H144").Value = "'4000.00"
$ws.Range("H145").Value = "'2500.00"
$ws.Range("H146").Value = "'27776.00"
$ws.Range("H147").Value = "'4000.00"
$ws.Range("H148").Value = "'1200.00"
$ws.Range("H149").Value = "'8400.00"
$ws.Range("H150").Value = "'11850.00"
$ws.Range("H151").Value = "'750.00"
$ws.Range("H152").Value = "'665.50"
$ws.Range("H153").Value = "'354.90"
$ws.Range("H154").Value = "'5269.00"
$ws.Range("H155").Value = "'36354.00"
$ws.Range("H156").Value = "'13160.00"
$ws.Range("H157").Value = "'351.46"
$ws.Range("H158").Value = "'286.60"
$ws.Range("H159").Value = "'10043.00"
$ws.Range("H160").Value = "'17500.00"
$ws.Range("H161").Value = "'30278.00"
$ws.Range("H162").Value = "'1000.00"
$ws.Range("H163").Value = "'16453.99"
$ws.Range("H164").Value = "'8927.00"
$ws.Range("H165").Value = "'9949.26"
$ws.Range("H166").Value = "'12630.00"
$ws.Range("H167").Value = "'2000.00"
$ws.Range("H168").Value = "'1157.28"
$ws.Range("H169").Value = "'5680.43"
$ws.Range("H170").Value = "'264.00"
$ws.Range("H171").Value = "'178415.00"
$ws.Range("H172").Value = "'70910.00"
$ws.Range("H173").Value = "'21156.20"
$ws.Range("H174").Value = "'19796.56"
$ws.Range("H175").Value = "'8500.00"
$ws.Range("H176").Value = "'25000.00"
$ws.Range("H177").Value = "'25000.00"
$ws.Range("H178").Value = "'25000.00"
$ws.Range("H179").Value = "'25000.00"
$ws.Range("H180").Value = "'25000.00"
$ws.Range("H181").Value = "'25000.00"
$ws.Range("H182").Value = "'50000.00"
$ws.Range("H183").Value = "'50000.00"
$ws.Range("H184").Value = "'50000.00"
$ws.Range("H185").Value = "'25000.00"
$ws.Range("H186").Value = "'3228646.85"
$ws.Range("H187").Value = "'28200.00"
$ws.Range("H188").Value = "'116000.00"
$ws.Range("H189").Value = "'1380.00"
$ws.Range("H190").Value = "'128000.00"
$ws.Range("H191").Value = "'137000.00"
$ws.Range("H192").Value = "'128000.00"
$ws.Range("H193").Value = "'128000.00"
$ws.Range("H194").Value = "'128000.00"
$ws.Range("H195").Value = "'128000.00"
$ws.Range("H196").Value = "'224000.00"
$ws.Range("H197").Value = "'224000.00"
$ws.Range("H198").Value = "'320000.00"
$ws.Range("H199").Value = "'128000.00"
$ws.Range("H200").Value = "'128000.00"
$ws.Range("H201").Value = "'128000.00"
$ws.Range("H202").Value = "'128000.00"
$ws.Range("H203").Value = "'128000.00"
$ws.Range("H204").Value = "'224000.00"
$ws.Range("H205").Value = "'320000.00"
$ws.Range("H206").Value = "'224000.00"
$ws.Range("H207").Value = "'128000.00"
$ws.Range("H208").Value = "'211000.00"
$ws.Range("H209").Value = "'128000.00"
$ws.Range("H210").Value = "'128000.00"
$ws.Range("H211").Value = "'132750.00"
$ws.Range("H212").Value = "'128000.00"
$ws.Range("H213").Value = "'344244.32"
$ws.Range("H214").Value = "'18500.00"
$ws.Range("H215").Value = "'9483.24"
$ws.Range("H216").Value = "'27468.00"
$ws.Range("H217").Value = "'17300.00"
$ws.Range("H218").Value = "'30000.00"
$ws.Range("H219").Value = "'2607.06"
$ws.Range("H220").Value = "'23205.00"
$ws.Range("H221").Value = "'425000.00"
$ws.Range("H222").Value = "'18453.50"
$ws.Range("H223").Value = "'350.00"
$ws.Range("H224").Value = "'47730.00"
$ws.Range("H225").Value = "'10800.00"
